$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels (bold, matching style of other header cells like D6/E6/D9/E9 etc.)
$ws.Range("D18").Value = "Mean increase"
$ws.Range("D18").Font.Bold = $true

$ws.Range("F18").Value = "Median increase"
$ws.Range("F18").Font.Bold = $true

# New formulas: percentage increase of the (new) mean / median vs. the older
# recorded baseline values.
$ws.Range("D19").Formula = "= ((E3 / 95.321842) * 100) - 100"
$ws.Range("D19").ClearFormats()

$ws.Range("F19").Formula = "= ((E10 / 95.22216) * 100) - 100"
$ws.Range("F19").ClearFormats()

# Leave the selection where the author's last edit ended up.
[void]$ws.Range("F20").Select()
